$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(31).Insert()

$ws.Range("A31").Value = 9
$ws.Range("B31").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C31").Value = "Metropolitana"
$ws.Range("D31").Value = 44690
$ws.Range("E31").Value = 13
$ws.Range("F31").Value = 100114007
$ws.Range("G31").Value = "Jengibre"
$ws.Range("H31").Value = "Sin especificar"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 790
$ws.Range("K31").Value = 12000
$ws.Range("L31").Value = 13000
$ws.Range("M31").Value = 12494
$ws.Range("N31").Value = "$/caja 13 kilos"
$ws.Range("O31").Value = "Perú"
$ws.Range("P31").Value = 961
$ws.Range("Q31").Value = 13
$ws.Range("R31").Value = "Hortaliza"
